$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.173.35'
$ws.Range("E2").Value = '  -3.04%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.863.31'
$ws.Range("E3").Value = '  -3.92%  '

$ws.Range("E4").Value = '  -0.10%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '234.01'
$ws.Range("E5").Value = '  -3.34%  '

$ws.Range("E6").Value = '  -0.09%  '

$ws.Range("E7").Value = '  -2.63%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2828'
$ws.Range("E8").Value = '  -2.86%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06542'
$ws.Range("E9").Value = '  -3.53%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.21'
$ws.Range("E10").Value = '  +0.47%  '

$ws.Range("E11").Value = '  -0.34%  '

$ws.Range("E12").Value = '  -7.90%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.854.91'
$ws.Range("E13").Value = '  -4.39%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.128'
$ws.Range("E14").Value = '  -3.27%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6704'
$ws.Range("E15").Value = '  -3.87%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '279.85'
$ws.Range("E16").Value = '  -5.39%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '30.204.84'
$ws.Range("E17").Value = '  -3.00%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.000'

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.458'
$ws.Range("E19").Value = '  -2.01%  '

$ws.Range("E20").Value = '  -2.65%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.098.95'
$ws.Range("E21").Value = '  -4.62%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.000007246'
$ws.Range("E22").Value = '  -4.65%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.001'
$ws.Range("E23").Value = '  -0.08%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.150'
$ws.Range("E24").Value = '  -4.11%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.320'
$ws.Range("E25").Value = '  -2.41%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '165.34'
$ws.Range("E26").Value = '  -2.18%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.91'
$ws.Range("E27").Value = '  -4.60%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.905'
$ws.Range("E28").Value = '  -9.08%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.343'
$ws.Range("E29").Value = '  -3.64%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.09595'
$ws.Range("E30").Value = '  -4.81%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.405'
$ws.Range("E31").Value = '  -4.70%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.470'
$ws.Range("E32").Value = '  -4.25%  '

$ws.Range("E33").Value = '  -5.08%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04668'
$ws.Range("E34").Value = '  -3.60%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7018'
$ws.Range("E35").Value = '  -4.95%  '

$ws.Range("E36").Value = '  -3.02%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.711'
$ws.Range("E37").Value = '  -0.50%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01870'
$ws.Range("E38").Value = '  -4.74%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.275'
$ws.Range("E39").Value = '  -8.64%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.530'
$ws.Range("E40").Value = '  -3.99%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '72.46'
$ws.Range("E41").Value = '  -5.09%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8522'
$ws.Range("E42").Value = '  -2.44%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.930'
$ws.Range("E43").Value = '  -5.18%  '

$ws.Range("B44").Value = 'PaxDollar'
$ws.Range("C44").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.0000'
$ws.Range("E44").Value = '  -0.04%  '

$ws.Range("B45").Value = 'TheSandbox'
$ws.Range("C45").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.4162'
$ws.Range("E45").Value = '  -4.69%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '103.20'
$ws.Range("E46").Value = '  -2.45%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '988.25'
$ws.Range("E47").Value = '  -4.07%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.135'
$ws.Range("E48").Value = '  -5.81%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.205'
$ws.Range("E49").Value = '  -0.48%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '34.03'
$ws.Range("E50").Value = '  -3.46%  '
